# Data Importer test workbook fixes/enhancements (#2817 / #2818)
#
# The sheet gains a new column "M" holding relative-property-path sample
# values; the previous column M ("path" / "/tmp/node1" / "/tmp/node2")
# shifts right to become column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at M. This shifts the old M column (and its data)
# one position to the right, turning it into column N.
$ws.Columns("M:M").Insert() | Out-Null

# Populate the newly inserted column M with the new sample data.
$ws.Range("M1").Value = "./foo/bar/test"
$ws.Range("M2").Value = "relative property path 1"
$ws.Range("M3").Value = "relative property path 2"

# The inserted column picks up formatting from its left neighbour (L) for
# the data rows; strip that back off so M2/M3 stay plain/unstyled, matching
# the header-only style (M1 keeps the bold header style already present).
$ws.Range("M2:M3").ClearFormats() | Out-Null

# Give the new column a sensible explicit width (matches the other
# "bestFit"-ish wide text columns in the sheet).
$ws.Columns("M:M").ColumnWidth = 24.83

# Select the whole (shifted) column N, as was left selected in the source
# workbook after the edit.
$ws.Columns("N:N").Select() | Out-Null
